# Update player standings after a newly added match.
# Totals (TOTAL RUNS in column B, WICKETS in column C) are incremented
# for the players who featured in the new match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 29.0
$ws.Range("B4").Value = 80.0
$ws.Range("B5").Value = 8.0
$ws.Range("C8").Value = 4.0
$ws.Range("B9").Value = 33.0
$ws.Range("B10").Value = 35.0
$ws.Range("C10").Value = 6.0
$ws.Range("B11").Value = 28.0
$ws.Range("B12").Value = 24.0
$ws.Range("C12").Value = 7.0
$ws.Range("B46").Value = 96.0
$ws.Range("B47").Value = 100.0
$ws.Range("B48").Value = 238.0
$ws.Range("B49").Value = 183.0
$ws.Range("B50").Value = 157.0
$ws.Range("B51").Value = 60.0
$ws.Range("B52").Value = 102.0
$ws.Range("B53").Value = 46.0
$ws.Range("B54").Value = 146.0
$ws.Range("C55").Value = 15.0
$ws.Range("B56").Value = 48.0
$ws.Range("C56").Value = 22.0
$ws.Range("B61").Value = 263.0
$ws.Range("B62").Value = 99.0
$ws.Range("C63").Value = 12.0
$ws.Range("C64").Value = 7.0
$ws.Range("B65").Value = 24.0
$ws.Range("C65").Value = 9.0
$ws.Range("B66").Value = 33.0
$ws.Range("C66").Value = 9.0
$ws.Range("C67").Value = 13.0
$ws.Range("B70").Value = 92.0
$ws.Range("B71").Value = 95.0
$ws.Range("B72").Value = 83.0
$ws.Range("B73").Value = 28.0
$ws.Range("C74").Value = 8.0
$ws.Range("B75").Value = 115.0
$ws.Range("C75").Value = 7.0
$ws.Range("C76").Value = 9.0
$ws.Range("C77").Value = 7.0
$ws.Range("C78").Value = 8.0
